$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix row 12's Result (column A) which was previously "Fail" and should be "Pass"
$ws.Range("A12").Value = "Pass"

# Update the Date column (column B) for rows 2-22 with new timestamps
$ws.Range("B2").Value = "Thu Feb 16 12:16:31 EST 2023"
$ws.Range("B3").Value = "Thu Feb 16 12:17:14 EST 2023"
$ws.Range("B4").Value = "Thu Feb 16 12:18:02 EST 2023"
$ws.Range("B5").Value = "Thu Feb 16 12:18:43 EST 2023"
$ws.Range("B6").Value = "Thu Feb 16 12:19:24 EST 2023"
$ws.Range("B7").Value = "Thu Feb 16 12:20:05 EST 2023"
$ws.Range("B8").Value = "Thu Feb 16 12:20:48 EST 2023"
$ws.Range("B9").Value = "Thu Feb 16 12:21:37 EST 2023"
$ws.Range("B10").Value = "Thu Feb 16 12:22:24 EST 2023"
$ws.Range("B11").Value = "Thu Feb 16 12:23:12 EST 2023"
$ws.Range("B12").Value = "Thu Feb 16 12:23:55 EST 2023"
$ws.Range("B13").Value = "Thu Feb 16 12:24:43 EST 2023"
$ws.Range("B14").Value = "Thu Feb 16 12:25:26 EST 2023"
$ws.Range("B15").Value = "Thu Feb 16 12:26:11 EST 2023"
$ws.Range("B16").Value = "Thu Feb 16 12:26:57 EST 2023"
$ws.Range("B17").Value = "Thu Feb 16 12:27:42 EST 2023"
$ws.Range("B18").Value = "Thu Feb 16 12:28:30 EST 2023"
$ws.Range("B19").Value = "Thu Feb 16 12:29:19 EST 2023"
$ws.Range("B20").Value = "Thu Feb 16 12:30:07 EST 2023"
$ws.Range("B21").Value = "Thu Feb 16 12:30:53 EST 2023"
$ws.Range("B22").Value = "Thu Feb 16 12:31:40 EST 2023"
